# Update countries & provincias Spain
#
# Refresh the COVID-19 "Pais" worksheet with newer case figures:
#  - bump the "Datos actualizados" timestamp (row 1)
#  - update Casos totales/Nuevos casos/Casos activos/Recuperados/Casos
#    criticos/Muertes hoy/Muertes (columns B-H) for the countries whose
#    figures changed
#  - because Indonesia, Eslovenia and Letonia overtake their neighbours
#    in the ranking, columns A-H are rewritten for the small block of
#    rows where the row order swaps (22-24, 112-114, 151-152)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: "last updated" timestamp
$ws.Range("A1").Value = "Datos actualizados a 13 de Octubre de 2020 a las 12:06"


# Row 19: Banglades
$ws.Range("B19").Value = 381275
$ws.Range("C19").Value = 1537
$ws.Range("D19").Value = 295873
$ws.Range("E19").Value = 79825
$ws.Range("F19").Value = 0
$ws.Range("G19").Value = 22
$ws.Range("H19").Value = 5577

# Row 22: Indonesia (was Arabia Saudita)
$ws.Range("A22").Value = "Indonesia"
$ws.Range("B22").Value = 340622
$ws.Range("C22").Value = 3906
$ws.Range("D22").Value = 263296
$ws.Range("E22").Value = 65299
$ws.Range("F22").Value = 0
$ws.Range("G22").Value = 92
$ws.Range("H22").Value = 12027

# Row 23: Arabia Saudita (was Turquia)
$ws.Range("A23").Value = "Arabia Saudita"
$ws.Range("B23").Value = 339615
$ws.Range("C23").Value = 0
$ws.Range("D23").Value = 325839
$ws.Range("E23").Value = 8708
$ws.Range("F23").Value = 0
$ws.Range("G23").Value = 0
$ws.Range("H23").Value = 5068

# Row 24: Turquia (was Indonesia)
$ws.Range("A24").Value = "Turquia"
$ws.Range("B24").Value = 337147
$ws.Range("C24").Value = 0
$ws.Range("D24").Value = 295658
$ws.Range("E24").Value = 32594
$ws.Range("F24").Value = 0
$ws.Range("G24").Value = 0
$ws.Range("H24").Value = 8895

# Row 64: Austria
$ws.Range("B64").Value = 57326
$ws.Range("C64").Value = 1028
$ws.Range("D64").Value = 44949
$ws.Range("E64").Value = 11516
$ws.Range("F64").Value = 0
$ws.Range("G64").Value = 6
$ws.Range("H64").Value = 861

# Row 83: Australia
$ws.Range("B83").Value = 27316
$ws.Range("C83").Value = 30
$ws.Range("D83").Value = 25037
$ws.Range("E83").Value = 1380
$ws.Range("F83").Value = 0
$ws.Range("G83").Value = 1
$ws.Range("H83").Value = 899

# Row 102: Finlandia
$ws.Range("B102").Value = 12499
$ws.Range("C102").Value = 287
$ws.Range("D102").Value = 8500
$ws.Range("E102").Value = 3653
$ws.Range("F102").Value = 0
$ws.Range("G102").Value = 0
$ws.Range("H102").Value = 346

# Row 104: Guinea
$ws.Range("B104").Value = 11134
$ws.Range("C104").Value = 72
$ws.Range("D104").Value = 10347
$ws.Range("E104").Value = 717
$ws.Range("F104").Value = 0
$ws.Range("G104").Value = 0
$ws.Range("H104").Value = 70

# Row 112: Eslovenia (was Haiti)
$ws.Range("A112").Value = "Eslovenia"
$ws.Range("B112").Value = 9231
$ws.Range("C112").Value = 398
$ws.Range("D112").Value = 5314
$ws.Range("E112").Value = 3744
$ws.Range("F112").Value = 0
$ws.Range("G112").Value = 4
$ws.Range("H112").Value = 173

# Row 113: Haiti (was Gabon)
$ws.Range("A113").Value = "Haiti"
$ws.Range("B113").Value = 8882
$ws.Range("C113").Value = 0
$ws.Range("D113").Value = 7104
$ws.Range("E113").Value = 1548
$ws.Range("F113").Value = 0
$ws.Range("G113").Value = 0
$ws.Range("H113").Value = 230

# Row 114: Gabon (was Eslovenia)
$ws.Range("A114").Value = "Gabon"
$ws.Range("B114").Value = 8860
$ws.Range("C114").Value = 0
$ws.Range("D114").Value = 8298
$ws.Range("E114").Value = 508
$ws.Range("F114").Value = 0
$ws.Range("G114").Value = 0
$ws.Range("H114").Value = 54

# Row 121: Lituania
$ws.Range("B121").Value = 6366
$ws.Range("C121").Value = 118
$ws.Range("D121").Value = 2842
$ws.Range("E121").Value = 3418
$ws.Range("F121").Value = 0
$ws.Range("G121").Value = 3
$ws.Range("H121").Value = 106

# Row 127: Hong Kong
$ws.Range("B127").Value = 5202
$ws.Range("C127").Value = 8
$ws.Range("D127").Value = 4921
$ws.Range("E127").Value = 176
$ws.Range("F127").Value = 0
$ws.Range("G127").Value = 0
$ws.Range("H127").Value = 105

# Row 151: Letonia (was Sudan del Sur)
$ws.Range("A151").Value = "Letonia"
$ws.Range("B151").Value = 2840
$ws.Range("C151").Value = 75
$ws.Range("D151").Value = 1325
$ws.Range("E151").Value = 1474
$ws.Range("F151").Value = 0
$ws.Range("G151").Value = 1
$ws.Range("H151").Value = 41

# Row 152: Sudan del Sur (was Letonia)
$ws.Range("A152").Value = "Sudan del Sur"
$ws.Range("B152").Value = 2787
$ws.Range("C152").Value = 0
$ws.Range("D152").Value = 1290
$ws.Range("E152").Value = 1442
$ws.Range("F152").Value = 0
$ws.Range("G152").Value = 0
$ws.Range("H152").Value = 55

# Row 176: Taiwan
$ws.Range("B176").Value = 530
$ws.Range("C176").Value = 1
$ws.Range("D176").Value = 489
$ws.Range("E176").Value = 34
$ws.Range("F176").Value = 0
$ws.Range("G176").Value = 0
$ws.Range("H176").Value = 7

# Row 185: Isla de Man
$ws.Range("B185").Value = 346
$ws.Range("C185").Value = 0
$ws.Range("D185").Value = 317
$ws.Range("E185").Value = 5
$ws.Range("F185").Value = 0
$ws.Range("G185").Value = 0
$ws.Range("H185").Value = 24

# Row 212: San Cristobal y Nieves
$ws.Range("B212").Value = 19
$ws.Range("C212").Value = 0
$ws.Range("D212").Value = 19
$ws.Range("E212").Value = 0
$ws.Range("F212").Value = 0
$ws.Range("G212").Value = 0
$ws.Range("H212").Value = 0
